# Append the latest portfolio snapshot row (2025-10-27) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 73

# Column A holds a plain "YYYY-MM-DD" text label (the sheet stores dates as
# literal strings, not Excel date serials). Force text entry with a "@"
# number format so Excel doesn't auto-convert the literal into a date
# serial, then drop the format back to the sheet's default (no explicit
# style) so the new cell matches the rest of the column.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2025-10-27"
$ws.Range("A" + $newRow).Style = "Normal"

# Columns B/C are plain numeric quotes for this date (SUZLON.NS / TATAMOTORS.NS).
# Column D (ETERNAL.NS) has no quote for this date, so it stays empty.
$ws.Range("B" + $newRow).Value = 53.70999908447266
$ws.Range("C" + $newRow).Value = 333.7000122070312
